# Weekly update: insert a new price record for "Vega Monumental Concepción -
# Zanahoria" as the new row 203, pushing all the existing rows (old 203-224)
# down by one (new 204-225).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 203, shifting rows 203:224 -> 204:225.
$ws.Rows.Item(203).Insert()

# Populate the new row 203 with the latest weekly record. Most columns
# (Mercado ID, Mercado, Región, Codreg, Categoría ID, Categoría, Variedad,
# Calidad, Unidad de comercialización, Kg o Unidades, Clasificación) repeat
# the values used throughout this sub-series; only the date, volume,
# min/max/avg prices, origin and $/Kg price are new for this entry.
$ws.Range("A203").Value = 11
$ws.Range("B203").Value = "Vega Monumental Concepción"
$ws.Range("C203").Value = "Bíobío"
$ws.Range("D203").Value = 44714
$ws.Range("E203").Value = 8
$ws.Range("F203").Value = 100114013
$ws.Range("G203").Value = "Zanahoria"
$ws.Range("H203").Value = "Sin especificar"
$ws.Range("I203").Value = "Primera"
$ws.Range("J203").Value = 350
$ws.Range("K203").Value = 5000
$ws.Range("L203").Value = 5500
$ws.Range("M203").Value = 5286
$ws.Range("N203").Value = "$/saco 20 kilos"
$ws.Range("O203").Value = "Chillán"
$ws.Range("P203").Value = 264
$ws.Range("Q203").Value = 20
$ws.Range("R203").Value = "Hortaliza"
